$wb = $excel.ActiveWorkbook

# --- Rename headers on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# --- Header row ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match header formatting (bold, centered, bordered) used by the other sheets
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# --- Data rows ---
$forecastRows = @(
  @(45109.99999999999, 309, -250.7916786480656, 913.2631355489196),
  @(45123.99999999999, 314, -273.1235718813712, 883.5850855484028),
  @(45130.99999999999, 317, -256.6704213090392, 870.5724970014344),
  @(45137.99999999999, 320, -258.181823854159, 948.8453748944132),
  @(45144.99999999999, 323, -233.3512532737852, 955.0719762720428),
  @(45151.99999999999, 325, -225.1385779795161, 959.6581855951936),
  @(45186.99999999999, 339, -205.3474384200671, 916.9359713745482),
  @(45193.99999999999, 342, -225.8464851456566, 910.5072721600735),
  @(45200.99999999999, 344, -238.4302598470239, 924.7826106186926),
  @(45207.99999999999, 347, -199.4865620999718, 914.8063922979445),
  @(45214.99999999999, 350, -220.0230208910144, 941.0388441139095),
  @(45221.99999999999, 352, -245.2641339340131, 933.9104791153678),
  @(45228.99999999999, 355, -225.8053398712694, 924.5919044817784),
  @(45242.99999999999, 361, -223.0207689043973, 944.4529508535428),
  @(45256.99999999999, 366, -250.8977055205335, 935.9568425101122),
  @(45270.99999999999, 372, -227.467874841961, 978.7842475049589),
  @(45277.99999999999, 374, -176.7185734113199, 962.3813339411577),
  @(45298.99999999999, 382, -220.2980109716253, 950.8638122189185),
  @(45312.99999999999, 388, -153.6929748344431, 965.967004641725),
  @(45319.99999999999, 391, -163.4218357951243, 993.1568626399242),
  @(45326.99999999999, 393, -181.4218664264652, 952.5716935487649),
  @(45361.99999999999, 407, -212.1368489762532, 1010.124754314806),
  @(45375.99999999999, 412, -176.649413015434, 968.6284332524298),
  @(45417.99999999999, 429, -109.4732887876597, 1052.237687128632),
  @(45438.99999999999, 437, -127.2075127690815, 1031.080037534932),
  @(45445.99999999999, 440, -134.7049052913397, 1031.035134892128),
  @(45452.99999999999, 442, -126.4458246979479, 1023.991578598071),
  @(45459.99999999999, 445, -135.3944679880381, 1032.660954593352),
  @(45466.99999999999, 448, -138.6598001540828, 1031.73693849971),
  @(45473.99999999999, 451, -100.4064164151576, 1040.578716791754),
  @(45480.99999999999, 453, -167.6956575697383, 1006.198544785055),
  @(45487.99999999999, 456, -118.5506667968104, 1029.529026813922),
  @(45494.99999999999, 459, -163.5743283315645, 1045.461078076515),
  @(45501.99999999999, 461, -121.9035286671905, 1044.145770372587),
  @(45515.99999999999, 467, -109.3804109152147, 1041.915410162231),
  @(45522.99999999999, 470, -130.516878348967, 1000.822996856299),
  @(45529.99999999999, 472, -137.1011057673529, 1039.801882503038),
  @(45543.99999999999, 478, -140.2211115433528, 1076.017623856826),
  @(45550.99999999999, 480, -69.30760120866564, 1021.109025662187),
  @(45557.99999999999, 483, -110.5352846254615, 1088.622549777875),
  @(45585.99999999999, 494, -60.117867883183, 1073.935965756042),
  @(45599.99999999999, 500, -65.83964830892228, 1105.251204131551),
  @(45606.99999999999, 502, -87.32944540958616, 1091.356632813293),
  @(45613.99999999999, 505, -90.45045181803685, 1101.681587985677),
  @(45620.99999999999, 508, -101.8967683821441, 1094.42185467911),
  @(45627.99999999999, 510, -94.67589912015116, 1073.366734199651),
  @(45634.99999999999, 513, -61.31505305322703, 1093.665537522543),
  @(45641.99999999999, 516, -48.51484782212163, 1088.683641697648),
  @(45648.99999999999, 519, -78.7959095900962, 1101.268441861977),
  @(45655.99999999999, 521, -55.62657217844387, 1097.30721430005)
)

$r = 2
foreach ($row in $forecastRows) {
  $wsForecast.Cells.Item($r, 1).Value = $row[0]
  $wsForecast.Cells.Item($r, 2).Value = $row[1]
  $wsForecast.Cells.Item($r, 3).Value = $row[2]
  $wsForecast.Cells.Item($r, 4).Value = $row[3]
  $r = $r + 1
}

# Match the date-formatted style used for the "ds"/date column elsewhere
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A51").PasteSpecial(-4122)

$null = $wsForecast.Range("A1").Select()
